$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 798
$ws.Range("L3").Value = 833
$ws.Range("L4").Value = 208
$ws.Range("L5").Value = 57
$ws.Range("L6").Value = 886
$ws.Range("L7").Value = 2782

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 11
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 170

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 22
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 23
$ws.Range("L5").Value = 5
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 93
$ws.Range("L8").Value = 170
$ws.Range("L11").Value = 42
$ws.Range("L15").Value = 20
$ws.Range("L18").Value = 21
$ws.Range("L19").Value = 92
$ws.Range("L20").Value = 75
$ws.Range("L23").Value = 28
$ws.Range("L30").Value = 17
$ws.Range("L33").Value = 118
$ws.Range("L35").Value = 4
$ws.Range("L37").Value = 90
$ws.Range("L42").Value = 91
$ws.Range("L47").Value = 23
$ws.Range("L48").Value = 47
$ws.Range("L53").Value = 36
$ws.Range("L55").Value = 27
$ws.Range("L56").Value = 2
$ws.Range("L63").Value = 12
$ws.Range("L64").Value = 24
$ws.Range("L67").Value = 100
$ws.Range("L72").Value = 11
$ws.Range("L75").Value = 11
$ws.Range("L78").Value = 40
$ws.Range("L79").Value = 79
$ws.Range("L83").Value = 56
$ws.Range("L84").Value = 26
$ws.Range("L85").Value = 138
$ws.Range("L86").Value = 19
$ws.Range("L88").Value = 43
$ws.Range("L90").Value = 24
$ws.Range("L94").Value = 34
$ws.Range("L95").Value = 41
$ws.Range("L96").Value = 24
$ws.Range("L97").Value = 36
$ws.Range("L99").Value = 44
$ws.Range("L100").Value = 3
$ws.Range("L101").Value = 2782

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 27
$ws.Range("L3").Value = 28
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 29
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 17
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 10
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 3

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 33
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 4

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L2").Value = 2
$ws.Range("L6").Value = 19

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 65
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 11

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L3").Value = 1
$ws.Range("L6").Value = 2
